# Append one new data row (row 98) to each of the 4 worksheets, mirroring
# the existing "time / 总长 / ID / 实际长度 / 和校验 / ..._DEC" record layout.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: MID_LFT_#1 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A98").Value = 45884.46582175926
$ws1.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B98").Value = "0x01,0x90"
$ws1.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws1.Range("D98").Value = "0x01,0x14"
$ws1.Range("E98").Value = "0x07"
$ws1.Range("F98").Value = 400
$ws1.Range("G98").Value = [double]"5.68631262647113e+23"
$ws1.Range("H98").Value = 276
$ws1.Range("I98").Value = 7

# --- Sheet 2: MID_LFT_#2 ---------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A98").Value = 45884.46582175926
$ws2.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B98").Value = "0x01,0x7c"
$ws2.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D98").Value = "0x01,0x24"
$ws2.Range("E98").Value = "0x19"
$ws2.Range("F98").Value = 380
$ws2.Range("G98").Value = [double]"5.68432987514711e+23"
$ws2.Range("H98").Value = 292
$ws2.Range("I98").Value = 25

# --- Sheet 3: MID_PLT_#1 ---------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A98").Value = 45884.46582175926
$ws3.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B98").Value = "0x00,0x6e"
$ws3.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D98").Value = "0x00,0x5D"
$ws3.Range("E98").Value = "0x15"
$ws3.Range("F98").Value = 110
$ws3.Range("G98").Value = [double]"5.68631262647113e+23"
$ws3.Range("H98").Value = 93
$ws3.Range("I98").Value = 15

# --- Sheet 4: MID_PLT_#2 ---------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A98").Value = 45884.46582175926
$ws4.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B98").Value = "0x00,0x82"
$ws4.Range("C98").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D98").Value = "0x00,0x74"
$ws4.Range("E98").Value = "0x9"
$ws4.Range("F98").Value = 130
$ws4.Range("G98").Value = [double]"5.68631262647113e+23"
$ws4.Range("H98").Value = 116
$ws4.Range("I98").Value = 9
